$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 97: add D97 (reuses existing shared string "Revision of C-DS-OS concepts")
$ws.Range("D97").Value = "Revision of C-DS-OS concepts"

# Row 98: new status entry for 11/11/2021
$ws.Range("A98").Value = "11/11/2021"
$ws.Range("B98").Value = "Attended the video basics ppt by thenew teammate"
$ws.Range("C98").Value = "referring the shared links"

# Row 99
$ws.Range("B99").Value = "discussed the doubts,installed yuv players "
$ws.Range("C99").Value = "Updated the notes"

# Row 100
$ws.Range("B100").Value = "Attended the LDD recap session"

# Update selection to match final cursor position
[void]$ws.Range("D99").Select()

Write-Output "done"
